# Update loading_percent results for the 380 kV case (Case_5_226)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 12.83636662940317
$ws.Range("C2").Value = 9.858893280965644
$ws.Range("D2").Value = 6.035034581594623
$ws.Range("E2").Value = 12.40552867458604
$ws.Range("F2").Value = 29.17808867001389
$ws.Range("K2").Value = 9.200439912298622
$ws.Range("L2").Value = 9.629963664221798
$ws.Range("O2").Value = 26.16319783465897

# Row 3
$ws.Range("B3").Value = 12.58634283144126
$ws.Range("C3").Value = 9.869088147762357
$ws.Range("D3").Value = 5.99939281838004
$ws.Range("E3").Value = 12.41836552109813
$ws.Range("F3").Value = 29.21694017154272
$ws.Range("K3").Value = 9.019088439655929
$ws.Range("L3").Value = 9.616304712339197
$ws.Range("O3").Value = 26.23728306690399

# Row 4
$ws.Range("B4").Value = 12.43274874063101
$ws.Range("C4").Value = 9.875864392211966
$ws.Range("D4").Value = 5.977044282430203
$ws.Range("E4").Value = 12.42856340197524
$ws.Range("F4").Value = 29.247974055141
$ws.Range("K4").Value = 8.907323940884734
$ws.Range("L4").Value = 9.609573857461271
$ws.Range("O4").Value = 26.28812220415171

# Row 5
$ws.Range("B5").Value = 12.37022585299068
$ws.Range("C5").Value = 9.878756095433459
$ws.Range("D5").Value = 5.967822980051971
$ws.Range("E5").Value = 12.43330146478711
$ws.Range("F5").Value = 29.26242262540693
$ws.Range("K5").Value = 8.861736647058095
$ws.Range("L5").Value = 9.607249651936876
$ws.Range("O5").Value = 26.31018238476176

# Row 6
$ws.Range("B6").Value = 12.35985067408344
$ws.Range("C6").Value = 9.879244144951898
$ws.Range("D6").Value = 5.966284977222403
$ws.Range("E6").Value = 12.43412338432649
$ws.Range("F6").Value = 29.26493052958671
$ws.Range("K6").Value = 8.854166186231465
$ws.Range("L6").Value = 9.606889068266895
$ws.Range("O6").Value = 26.31392648501451

# Row 7
$ws.Range("B7").Value = 12.43190514288592
$ws.Range("C7").Value = 9.875902862456993
$ws.Range("D7").Value = 5.976920379149067
$ws.Range("E7").Value = 12.42862494342738
$ws.Range("F7").Value = 29.24816162209678
$ws.Range("K7").Value = 8.906709223536888
$ws.Range("L7").Value = 9.60954081437453
$ws.Range("O7").Value = 26.2884142820093

# Row 8
$ws.Range("B8").Value = 12.7502267473879
$ws.Range("C8").Value = 9.862301521174615
$ws.Range("D8").Value = 6.022843271790872
$ws.Range("E8").Value = 12.40947415838365
$ws.Range("F8").Value = 29.18999297856066
$ws.Range("K8").Value = 9.138032400014756
$ws.Range("L8").Value = 9.624911861240253
$ws.Range("O8").Value = 26.18763039704669

# Row 9
$ws.Range("B9").Value = 13.37008440802253
$ws.Range("C9").Value = 9.839708834549315
$ws.Range("D9").Value = 6.109098681552312
$ws.Range("E9").Value = 12.39029198177722
$ws.Range("F9").Value = 29.13299950006133
$ws.Range("K9").Value = 9.585746303022576
$ws.Range("L9").Value = 9.668080173699563
$ws.Range("O9").Value = 26.03255931855447

# Row 10
$ws.Range("B10").Value = 13.81800485084971
$ws.Range("C10").Value = 9.825571311995342
$ws.Range("D10").Value = 6.170004330189186
$ws.Range("E10").Value = 12.38738541562169
$ws.Range("F10").Value = 29.12603353050408
$ws.Range("K10").Value = 9.907714257184216
$ws.Range("L10").Value = 9.707570611228943
$ws.Range("O10").Value = 25.94472382462172

# Row 11
$ws.Range("B11").Value = 14.0191957713567
$ws.Range("C11").Value = 9.819668870438633
$ws.Range("D11").Value = 6.197145035164279
$ws.Range("E11").Value = 12.3884855831691
$ws.Range("F11").Value = 29.13045217163483
$ws.Range("K11").Value = 10.05201181297792
$ws.Range("L11").Value = 9.727183058447173
$ws.Range("O11").Value = 25.91045964479191

# Row 12
$ws.Range("B12").Value = 14.09493783731448
$ws.Range("C12").Value = 9.817509380484351
$ws.Range("D12").Value = 6.207338519019174
$ws.Range("E12").Value = 12.38924960488999
$ws.Range("F12").Value = 29.13321563444157
$ws.Range("K12").Value = 10.10629071537204
$ws.Range("L12").Value = 9.73484284916527
$ws.Range("O12").Value = 25.89830531072852

# Row 13
$ws.Range("B13").Value = 14.07864634749866
$ws.Range("C13").Value = 9.817971107695424
$ws.Range("D13").Value = 6.205146954713183
$ws.Range("E13").Value = 12.38906962471989
$ws.Range("F13").Value = 29.13257200826777
$ws.Range("K13").Value = 10.09461773837176
$ws.Range("L13").Value = 9.733182880733304
$ws.Range("O13").Value = 25.90088642763078

# Row 14
$ws.Range("B14").Value = 14.02543644034151
$ws.Range("C14").Value = 9.819489694322465
$ws.Range("D14").Value = 6.197985360758789
$ws.Range("E14").Value = 12.3885414830357
$ws.Range("F14").Value = 29.13065768172536
$ws.Range("K14").Value = 10.05648494053621
$ws.Range("L14").Value = 9.727808588578252
$ws.Range("O14").Value = 25.90944323925065

# Row 15
$ws.Range("B15").Value = 13.99278378964651
$ws.Range("C15").Value = 9.820429710714283
$ws.Range("D15").Value = 6.193587641825602
$ws.Range("E15").Value = 12.38826319351888
$ws.Range("F15").Value = 29.12962704010594
$ws.Range("K15").Value = 10.0330786852031
$ws.Range("L15").Value = 9.724546897391814
$ws.Range("O15").Value = 25.9147914827932

# Row 16
$ws.Range("B16").Value = 13.80479832256563
$ws.Range("C16").Value = 9.825967656502122
$ws.Range("D16").Value = 6.168219001733175
$ws.Range("E16").Value = 12.38736218720234
$ws.Range("F16").Value = 29.12589738847683
$ws.Range("K16").Value = 9.898235962035013
$ws.Range("L16").Value = 9.706321718239138
$ws.Range("O16").Value = 25.94707780981774

# Row 17
$ws.Range("B17").Value = 13.68876229890597
$ws.Range("C17").Value = 9.829500154551237
$ws.Range("D17").Value = 6.152509313564454
$ws.Range("E17").Value = 12.38742928499015
$ws.Range("F17").Value = 29.12555237403592
$ws.Range("K17").Value = 9.814921382407952
$ws.Range("L17").Value = 9.695560499610556
$ws.Range("O17").Value = 25.96834405177647

# Row 18
$ws.Range("B18").Value = 13.62178372926203
$ws.Range("C18").Value = 9.831581749868249
$ws.Range("D18").Value = 6.143420445675921
$ws.Range("E18").Value = 12.38769595149346
$ws.Range("F18").Value = 29.12606820726949
$ws.Range("K18").Value = 9.766799947875047
$ws.Range("L18").Value = 9.689526309928615
$ws.Range("O18").Value = 25.98111153046172

# Row 19
$ws.Range("B19").Value = 13.5990674630814
$ws.Range("C19").Value = 9.832295107097915
$ws.Range("D19").Value = 6.140334081574568
$ws.Range("E19").Value = 12.38782543691023
$ws.Range("F19").Value = 29.12636554911679
$ws.Range("K19").Value = 9.750473948707349
$ws.Range("L19").Value = 9.687510041167007
$ws.Range("O19").Value = 25.98552631679577

# Row 20
$ws.Range("B20").Value = 13.70113969140316
$ws.Range("C20").Value = 9.829118963345403
$ws.Range("D20").Value = 6.154187149570756
$ws.Range("E20").Value = 12.38739854360695
$ws.Range("F20").Value = 29.12551517990276
$ws.Range("K20").Value = 9.823811559052258
$ws.Range("L20").Value = 9.696689996785434
$ws.Range("O20").Value = 25.96602476481807

# Row 21
$ws.Range("B21").Value = 14.04107812107455
$ws.Range("C21").Value = 9.819041599052937
$ws.Range("D21").Value = 6.20009119927723
$ws.Range("E21").Value = 12.38868719073185
$ws.Range("F21").Value = 29.13119038867985
$ws.Range("K21").Value = 10.0676957287491
$ws.Range("L21").Value = 9.729380858575315
$ws.Range("O21").Value = 25.90690760450027

# Row 22
$ws.Range("B22").Value = 14.26062459286743
$ws.Range("C22").Value = 9.812896142730313
$ws.Range("D22").Value = 6.229600615043286
$ws.Range("E22").Value = 12.39155381044278
$ws.Range("F22").Value = 29.14125309031447
$ws.Range("K22").Value = 10.22494787480105
$ws.Range("L22").Value = 9.752102189512168
$ws.Range("O22").Value = 25.87305605275973

# Row 23
$ws.Range("B23").Value = 14.14371213556319
$ws.Range("C23").Value = 9.816135897803175
$ws.Range("D23").Value = 6.213896769724544
$ws.Range("E23").Value = 12.38983896413579
$ws.Range("F23").Value = 29.13530158566873
$ws.Range("K23").Value = 10.14123157447198
$ws.Range("L23").Value = 9.739852706195697
$ws.Range("O23").Value = 25.89068477361475

# Row 24
$ws.Range("B24").Value = 13.69554470099799
$ws.Range("C24").Value = 9.829291141883925
$ws.Range("D24").Value = 6.153428777730431
$ws.Range("E24").Value = 12.38741173124379
$ws.Range("F24").Value = 29.12552977065028
$ws.Range("K24").Value = 9.819793000106616
$ws.Range("L24").Value = 9.696178875552294
$ws.Range("O24").Value = 25.96707162868273

# Row 25
$ws.Range("B25").Value = 13.20337514140252
$ws.Range("C25").Value = 9.845386682832755
$ws.Range("D25").Value = 6.086186615065082
$ws.Range("E25").Value = 12.39351505479327
$ws.Range("F25").Value = 29.14229270086267
$ws.Range("K25").Value = 9.465617528423008
$ws.Range("L25").Value = 9.655023561689378
$ws.Range("O25").Value = 26.06993748136787

